$d = $word.ActiveDocument

$d.Content.Find.Execute("92÷6=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "84÷8=10, 4", 2) | Out-Null
$d.Content.Find.Execute("18÷2=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷9=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=11, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷9=2, 8", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("47÷5=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=10, 5", 2) | Out-Null
$d.Content.Find.Execute("70÷7=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "95÷6=15, 5", 2) | Out-Null
$d.Content.Find.Execute("97÷7=13, 6", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=16, 0", 2) | Out-Null
$d.Content.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "81÷9=9, 0", 2) | Out-Null
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=48, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷5=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "89÷3=29, 2", 2) | Out-Null
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "75÷2=37, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷3=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷6=1, 5", 2) | Out-Null
$d.Content.Find.Execute("53÷7=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "23÷6=3, 5", 2) | Out-Null
$d.Content.Find.Execute("75÷5=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=3, 3", 2) | Out-Null
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=7, 0", 2) | Out-Null
$d.Content.Find.Execute("53÷4=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "91÷5=18, 1", 2) | Out-Null
$d.Content.Find.Execute("30÷2=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷6=11, 1", 2) | Out-Null
$d.Content.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=14, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷4=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷4=19, 3", 2) | Out-Null
$d.Content.Find.Execute("75÷3=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=17, 3", 2) | Out-Null
$d.Content.Find.Execute("19÷3=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=10, 0", 2) | Out-Null
$d.Content.Find.Execute("74÷3=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "95÷9=10, 5", 2) | Out-Null
$d.Content.Find.Execute("91÷8=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=5, 1", 2) | Out-Null
$d.Content.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "13÷8=1, 5", 2) | Out-Null
$d.Content.Find.Execute("50÷3=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
